# This script reproduces an updated cryptocurrency price/volume snapshot
# (new "Price" values in column D, new "Volume(1h)" values in column E)
# for rows 2-51 of the active worksheet.
#
# Column D values are written with a leading apostrophe. That is the standard
# Excel convention for forcing a numeric-looking entry ("212.01", "7.40", ...) to
# be stored as literal text instead of being parsed into a floating point number
# (which would silently drop trailing zeros / significant digits, e.g. turn
# "7.40" into 7.4). The apostrophe itself is not stored as part of the cell text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''28.388.81'
$ws.Range('E2').Value = '  +2.57%  '
$ws.Range('D3').Value = '''1.579.37'
$ws.Range('E3').Value = '  +0.41%  '
$ws.Range('E4').Value = '  +1.57%  '
$ws.Range('D5').Value = '''212.01'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('E6').Value = '  -0.31%  '
$ws.Range('E7').Value = '  +1.39%  '
$ws.Range('D8').Value = '''46.18'
$ws.Range('E8').Value = '  +3.13%  '
$ws.Range('D9').Value = '''24.04'
$ws.Range('E9').Value = '  +3.26%  '
$ws.Range('D10').Value = '''0.249'
$ws.Range('E10').Value = '  -0.68%  '
$ws.Range('D11').Value = '''0.0596'
$ws.Range('E11').Value = '  -0.26%  '
$ws.Range('D12').Value = '''0.0882'
$ws.Range('E12').Value = '  +0.84%  '
$ws.Range('D13').Value = '''1.805.22'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('D14').Value = '''1.567.67'
$ws.Range('E14').Value = '  -1.03%  '
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').Value = '''3.73'
$ws.Range('E16').Value = '  -1.10%  '
$ws.Range('D17').Value = '''28.422.40'
$ws.Range('E17').Value = '  +2.59%  '
$ws.Range('D18').Value = '''62.68'
$ws.Range('E18').Value = '  -0.78%  '
$ws.Range('D19').Value = '''229.08'
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').Value = '''0.0₃0700'
$ws.Range('E20').Value = '  -0.78%  '
$ws.Range('D21').Value = '''7.40'
$ws.Range('E21').Value = '  -1.51%  '
$ws.Range('E22').Value = '  +1.26%  '
$ws.Range('D23').Value = '''3.95'
$ws.Range('E23').Value = '  -4.13%  '
$ws.Range('D24').Value = '''9.24'
$ws.Range('E24').Value = '  -2.13%  '
$ws.Range('E25').Value = '  +3.10%  '
$ws.Range('D26').Value = '''151.11'
$ws.Range('E26').Value = '  +0.86%  '
$ws.Range('D27').Value = '''15.10'
$ws.Range('E27').Value = '  -1.41%  '
$ws.Range('D28').Value = '''6.51'
$ws.Range('E28').Value = '  -1.34%  '
$ws.Range('E29').Value = '  -1.85%  '
$ws.Range('E30').Value = '  +1.16%  '
$ws.Range('E31').Value = '  -1.43%  '
$ws.Range('D32').Value = '''0.0466'
$ws.Range('E32').Value = '  -1.51%  '
$ws.Range('E33').Value = '  -0.59%  '
$ws.Range('D34').Value = '''3.14'
$ws.Range('E34').Value = '  +0.43%  '
$ws.Range('D35').Value = '''1.390.76'
$ws.Range('E35').Value = '  -4.47%  '
$ws.Range('D36').Value = '''1.57'
$ws.Range('E36').Value = '  -1.43%  '
$ws.Range('E37').Value = '  -3.36%  '
$ws.Range('E38').Value = '  +1.21%  '
$ws.Range('E39').Value = '  +5.91%  '
$ws.Range('E40').Value = '  -1.09%  '
$ws.Range('D41').Value = '''0.539'
$ws.Range('E41').Value = '  -0.65%  '
$ws.Range('D42').Value = '''0.803'
$ws.Range('E42').Value = '  -1.25%  '
$ws.Range('E43').Value = '  +1.20%  '
$ws.Range('E44').Value = '  +0.80%  '
$ws.Range('E45').Value = '  -1.27%  '
$ws.Range('D46').Value = '''0.983'
$ws.Range('E46').Value = '  +1.16%  '
$ws.Range('D47').Value = '''62.70'
$ws.Range('E47').Value = '  -2.37%  '
$ws.Range('D48').Value = '''1.715.48'
$ws.Range('E48').Value = '  +0.61%  '
$ws.Range('D49').Value = '''86.18'
$ws.Range('E49').Value = '  -0.87%  '
$ws.Range('D50').Value = '''0.0₆0104'
$ws.Range('E50').Value = '  +4.42%  '
$ws.Range('D51').Value = '''0.0520'
$ws.Range('E51').Value = '  -0.67%  '
